$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '53.738.49'
$ws.Range('E2').Value = '  -8.71%  '
$ws.Range('D3').Value = '2.399.89'
$ws.Range('E3').Value = '  -12.15%  '
$ws.Range('E4').Value = '  +0.08%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '461.00'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -8.96%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '130.56'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -7.96%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.996'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -0.11%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.487'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -8.46%  '
$ws.Range('D9').Value = '2.415.28'
$ws.Range('E9').Value = '  -12.04%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0940'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -10.19%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '5.30'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -12.52%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.311'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -10.64%  '
$ws.Range('E13').Value = '  -4.17%  '
$ws.Range('D14').Value = '2.824.16'
$ws.Range('E14').Value = '  -12.07%  '
$ws.Range('D15').Value = '53.645.51'
$ws.Range('E15').Value = '  -9.02%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '19.57'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -9.97%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.0000131'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -3.87%  '
$ws.Range('D18').Value = '2.425.70'
$ws.Range('E18').Value = '  -11.26%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '4.14'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -13.00%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '305.53'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -11.29%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '9.34'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -15.35%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.991'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.78%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '5.67'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +1.14%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '5.30'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -15.30%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '55.69'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -12.12%  '
$ws.Range('E26').Value = '  +1.23%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '0.382'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -10.72%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.510.61'
$ws.Range('E28').Value = '  -11.96%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.153'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -11.27%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '7.00'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -7.02%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').Value = '0.0₃0708'
$ws.Range('E32').Value = '  -15.10%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '145.57'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -4.10%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '17.59'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -8.51%  '
$ws.Range('E35').Value = '  -11.69%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '4.93'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -9.06%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '3.48'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -16.94%  '
$ws.Range('E38').Value = '  -7.44%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.792'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -16.93%  '
$ws.Range('E40').Value = '  -0.22%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '32.65'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -9.40%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.590'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -2.42%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.0519'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -6.81%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '3.22'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -9.08%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '10.10'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -2.61%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.22'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -12.02%  '
$ws.Range('D47').Value = '1.922.42'
$ws.Range('E47').Value = '  -12.26%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0867'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -2.20%  '
$ws.Range('E49').Value = '  -4.70%  '
$ws.Range('E50').Value = '  -11.52%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '16.35'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -14.27%  '
